$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet/tab
$ws.Name = "CubeA"

# Add the new averaged-intensities row (Gaussian Quadrature Scheme result)
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"

# Copy the formatting (border/bold/alignment) used by the rest of column A
$ws.Range("A15").Copy() | Out-Null
$ws.Range("A16").PasteSpecial(-4122) | Out-Null

$ws.Range("C16").Value = 0.6227848949555937
$ws.Range("D16").Value = 1.642271734218754
$ws.Range("E16").Value = 0.9957675914803306
$ws.Range("F16").Value = 1.038766456686441
$ws.Range("G16").Value = 0.6227848949555937
$ws.Range("H16").Value = 1.642271734218754
$ws.Range("I16").Value = 0.8786467324694143
$ws.Range("J16").Value = 1.096356904225475
$ws.Range("K16").Value = 0.8079618251480029
$ws.Range("L16").Value = 1.313100730709999
$ws.Range("M16").Value = 0.6227848949555937
$ws.Range("N16").Value = 1.319019662849542
$ws.Range("O16").Value = 1.07489766933528
$ws.Range("P16").Value = 1.049457108736751
